$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 576 (shifts existing rows 576:605 down to 577:606,
# matching the diff which grows the sheet from A1:T605 to A1:T606).
$ws.Rows.Item(576).Insert()

# Populate the newly inserted row 576 with the new weekly record.
$ws.Range("A576").Value = 11
$ws.Range("B576").Value = 'Vega Monumental Concepción'
$ws.Range("C576").Value = 'Bíobío'
$ws.Range("D576").Value = 45267
$ws.Range("E576").Value = 8
$ws.Range("F576").Value = 'Fruta'
$ws.Range("G576").Value = 100101
$ws.Range("H576").Value = 'Berries'
$ws.Range("I576").Value = 100112025
$ws.Range("J576").Value = 'Frutilla'
$ws.Range("K576").Value = 'Sin especificar'
$ws.Range("L576").Value = 'Primera'
$ws.Range("M576").Value = 200
$ws.Range("N576").Value = 12000
$ws.Range("O576").Value = 12000
$ws.Range("P576").Value = 12000
$ws.Range("Q576").Value = '$/bandeja 7 kilos'
$ws.Range("R576").Value = 'Provincia de Melipilla'
$ws.Range("S576").Value = 1714
$ws.Range("T576").Value = 7
